$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Block($startRow, $endRow, $order) {
    $numRows = $endRow - $startRow + 1
    $rangeAddr = "F" + $startRow + ":V" + $endRow
    $orig = $ws.Range($rangeAddr).Value2
    for ($i = 0; $i -lt $numRows; $i++) {
        $targetRow = $startRow + $i
        $srcOffset = $order[$i]
        $arr = New-Object 'object[,]' 1,17
        for ($c = 1; $c -le 17; $c++) {
            $arr[0, $c-1] = $orig[$srcOffset+1, $c]
        }
        $targetAddr = "F" + $targetRow + ":V" + $targetRow
        $ws.Range($targetAddr).Value = $arr
    }
}

# Re-order existing rows (odds/results data was re-scraped and matches realigned)
Swap-Block 10 13 @(1,2,3,0)
Swap-Block 18 19 @(1,0)
Swap-Block 33 34 @(1,0)
Swap-Block 40 41 @(1,0)
Swap-Block 59 64 @(1,0,5,4,3,2)

# Append newly scraped matches (rows 76-82)

# Row 76
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = 75
$arr[0,1] = "netherlands"
$arr[0,2] = "tweede-divisie"
$arr[0,3] = "2023-2024"
$arr[0,4] = 45227.60416666666
$arr[0,5] = 'GVVV'
$arr[0,6] = 2
$arr[0,7] = 'Hardenberg'
$arr[0,8] = 1
$arr[0,9] = 2.57
$arr[0,10] = '27/10/2023 02:42'
$arr[0,11] = 2.53
$arr[0,12] = '28/10/2023 14:19'
$arr[0,13] = 3.4
$arr[0,14] = '27/10/2023 02:42'
$arr[0,15] = 3.58
$arr[0,16] = '28/10/2023 14:19'
$arr[0,17] = 2.25
$arr[0,18] = '27/10/2023 02:42'
$arr[0,19] = 2.44
$arr[0,20] = '28/10/2023 14:19'
$arr[0,21] = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/gvvv-hardenberg/Q3qCDdTk/'
$ws.Range("A75:V75").Copy()
$ws.Range("A76:V76").PasteSpecial(-4122)
$ws.Range("A76:V76").Value = $arr

# Row 77
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = 76
$arr[0,1] = "netherlands"
$arr[0,2] = "tweede-divisie"
$arr[0,3] = "2023-2024"
$arr[0,4] = 45227.60416666666
$arr[0,5] = 'Scheveningen'
$arr[0,6] = 0
$arr[0,7] = 'Katwijk'
$arr[0,8] = 4
$arr[0,9] = 3.32
$arr[0,10] = '27/10/2023 02:42'
$arr[0,11] = 2.94
$arr[0,12] = '28/10/2023 14:26'
$arr[0,13] = 3.55
$arr[0,14] = '27/10/2023 02:42'
$arr[0,15] = 3.74
$arr[0,16] = '28/10/2023 14:26'
$arr[0,17] = 1.83
$arr[0,18] = '27/10/2023 02:42'
$arr[0,19] = 2.1
$arr[0,20] = '28/10/2023 14:26'
$arr[0,21] = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/svv-scheveningen-katwijk/dSwyMKL9/'
$ws.Range("A75:V75").Copy()
$ws.Range("A77:V77").PasteSpecial(-4122)
$ws.Range("A77:V77").Value = $arr

# Row 78
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = 77
$arr[0,1] = "netherlands"
$arr[0,2] = "tweede-divisie"
$arr[0,3] = "2023-2024"
$arr[0,4] = 45227.625
$arr[0,5] = 'Kozakken Boys'
$arr[0,6] = 1
$arr[0,7] = 'Noordwijk'
$arr[0,8] = 3
$arr[0,9] = 2.26
$arr[0,10] = '27/10/2023 03:12'
$arr[0,11] = 2.49
$arr[0,12] = '28/10/2023 14:58'
$arr[0,13] = 3.34
$arr[0,14] = '27/10/2023 03:12'
$arr[0,15] = 3.47
$arr[0,16] = '28/10/2023 14:58'
$arr[0,17] = 2.6
$arr[0,18] = '27/10/2023 03:12'
$arr[0,19] = 2.53
$arr[0,20] = '28/10/2023 14:58'
$arr[0,21] = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/kozakken-boys-vv-noordwijk/4hsqKbjM/'
$ws.Range("A75:V75").Copy()
$ws.Range("A78:V78").PasteSpecial(-4122)
$ws.Range("A78:V78").Value = $arr

# Row 79
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = 78
$arr[0,1] = "netherlands"
$arr[0,2] = "tweede-divisie"
$arr[0,3] = "2023-2024"
$arr[0,4] = 45227.625
$arr[0,5] = 'Quick Boys'
$arr[0,6] = 3
$arr[0,7] = 'Lisse'
$arr[0,8] = 0
$arr[0,9] = 1.18
$arr[0,10] = '27/10/2023 03:12'
$arr[0,11] = 1.19
$arr[0,12] = '28/10/2023 14:59'
$arr[0,13] = 6.75
$arr[0,14] = '27/10/2023 03:12'
$arr[0,15] = 7.41
$arr[0,16] = '28/10/2023 14:59'
$arr[0,17] = 7.56
$arr[0,18] = '27/10/2023 03:12'
$arr[0,19] = 8.92
$arr[0,20] = '28/10/2023 14:59'
$arr[0,21] = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/quick-boys-lisse/nDp8ExEq/'
$ws.Range("A75:V75").Copy()
$ws.Range("A79:V79").PasteSpecial(-4122)
$ws.Range("A79:V79").Value = $arr

# Row 80
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = 79
$arr[0,1] = "netherlands"
$arr[0,2] = "tweede-divisie"
$arr[0,3] = "2023-2024"
$arr[0,4] = 45227.625
$arr[0,5] = 'Spakenburg'
$arr[0,6] = 4
$arr[0,7] = 'Jong Almere City'
$arr[0,8] = 3
$arr[0,9] = 1.46
$arr[0,10] = '27/10/2023 03:12'
$arr[0,11] = 1.43
$arr[0,12] = '28/10/2023 09:52'
$arr[0,13] = 4.52
$arr[0,14] = '27/10/2023 03:12'
$arr[0,15] = 4.97
$arr[0,16] = '28/10/2023 13:03'
$arr[0,17] = 4.5
$arr[0,18] = '27/10/2023 03:12'
$arr[0,19] = 5.33
$arr[0,20] = '28/10/2023 09:52'
$arr[0,21] = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/spakenburg-jong-almere-city/rqtmJI5S/'
$ws.Range("A75:V75").Copy()
$ws.Range("A80:V80").PasteSpecial(-4122)
$ws.Range("A80:V80").Value = $arr

# Row 81
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = 80
$arr[0,1] = "netherlands"
$arr[0,2] = "tweede-divisie"
$arr[0,3] = "2023-2024"
$arr[0,4] = 45227.64583333334
$arr[0,5] = 'Jong Sparta Rotterdam'
$arr[0,6] = 4
$arr[0,7] = 'HFC'
$arr[0,8] = 1
$arr[0,9] = 2.44
$arr[0,10] = '27/10/2023 03:43'
$arr[0,11] = 2.59
$arr[0,12] = '28/10/2023 15:10'
$arr[0,13] = 3.51
$arr[0,14] = '27/10/2023 03:43'
$arr[0,15] = 3.63
$arr[0,16] = '28/10/2023 15:10'
$arr[0,17] = 2.32
$arr[0,18] = '27/10/2023 03:43'
$arr[0,19] = 2.36
$arr[0,20] = '28/10/2023 15:10'
$arr[0,21] = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/jong-sparta-rotterdam-koninklijke-hfc/AVsuLvyG/'
$ws.Range("A75:V75").Copy()
$ws.Range("A81:V81").PasteSpecial(-4122)
$ws.Range("A81:V81").Value = $arr

# Row 82
$arr = New-Object 'object[,]' 1,22
$arr[0,0] = 81
$arr[0,1] = "netherlands"
$arr[0,2] = "tweede-divisie"
$arr[0,3] = "2023-2024"
$arr[0,4] = 45227.75
$arr[0,5] = 'De Treffers'
$arr[0,6] = 4
$arr[0,7] = 'Excelsior Maassluis'
$arr[0,8] = 0
$arr[0,9] = 1.34
$arr[0,10] = '27/10/2023 06:14'
$arr[0,11] = 1.33
$arr[0,12] = '28/10/2023 17:59'
$arr[0,13] = 5.2
$arr[0,14] = '27/10/2023 06:14'
$arr[0,15] = 5.51
$arr[0,16] = '28/10/2023 17:59'
$arr[0,17] = 5.51
$arr[0,18] = '27/10/2023 06:14'
$arr[0,19] = 6.5
$arr[0,20] = '28/10/2023 17:59'
$arr[0,21] = 'https://www.betexplorer.com/football/netherlands/tweede-divisie/de-treffers-excelsior-maassluis/YTnKBzb2/'
$ws.Range("A75:V75").Copy()
$ws.Range("A82:V82").PasteSpecial(-4122)
$ws.Range("A82:V82").Value = $arr
